$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear D3 value (previously 199.6231611810204) -> becomes a blank cell
$ws.Range("D3").ClearContents()

# Row 7: rename "Other" -> "Biogas" and update its value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 543.2088234238689

# New row 8: "Other" with its value, matching format/layout of row 7 (A3:A7 use style index 1)
$ws.Range("A7:D7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 254.3170937765217
